# Actualización automática 2025-08-07 16:50:08
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: VENTAS POR GRUPO ----
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M12").Value = 1090.97
$wsGrupo.Range("M24").Value = "1 de 22"

# ---- Sheet 2: VENTA MENSUAL ----
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F12").Value = 1090.97
$wsMensual.Range("F24").Value = 1090.97
# column F width 12 -> 13
$wsMensual.Columns.Item(6).ColumnWidth = 12.15

# ---- Sheet 3: CUMPLIMIENTO MENSUAL ----
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D16").Value = 1090.97
$wsCumpl.Range("E16").Value = 34965.73
$wsCumpl.Range("F16").Value = 0.0302570673411599

$wsCumpl.Range("D19").Value = 1090.97
$wsCumpl.Range("E19").Value = 53932.19386304604
$wsCumpl.Range("F19").Value = 0.0198274676228261

# column D width 11 -> 13
$wsCumpl.Columns.Item(4).ColumnWidth = 12.15
# column F width 18 -> 24
$wsCumpl.Columns.Item(6).ColumnWidth = 23.15
